$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 772778.1
$ws.Range("J17").Value = 772778.1
$ws.Range("L17").Value = 2318334.3
$ws.Range("N17").Value = -2318670.3
$ws.Range("H41").Value = 1267.1578
$ws.Range("J41").Value = 990.5833
$ws.Range("L41").Value = 990.5833
$ws.Range("N41").Value = -1870.5833
$ws.Range("H51").Value = 6747.6665
$ws.Range("J51").Value = 5998.273
$ws.Range("L51").Value = 5998.273
$ws.Range("N51").Value = -6966.273
$ws.Range("H52").Value = 2500
$ws.Range("J52").Value = 2500
$ws.Range("L52").Value = 7500
$ws.Range("N52").Value = -7820
$ws.Range("H98").Value = 1621
$ws.Range("I98").Value = 1647.8462
$ws.Range("J98").Value = 1446.5
$ws.Range("K98").Value = 1647.8462
$ws.Range("L98").Value = 1446.5
$ws.Range("M98").Value = -149.8462
$ws.Range("N98").Value = -4442.5
$ws.Range("H122").Value = 1621
$ws.Range("I122").Value = 1647.8462
$ws.Range("J122").Value = 1446.5
$ws.Range("K122").Value = 4943.5386
$ws.Range("L122").Value = 4339.5
$ws.Range("M122").Value = -2493.5386
$ws.Range("N122").Value = -9239.5
$ws.Range("H137").Value = 2192.4187
$ws.Range("I137").Value = 1935.8438
$ws.Range("J137").Value = 2938.818
$ws.Range("K137").Value = 5807.5314
$ws.Range("L137").Value = 8816.454000000002
$ws.Range("M137").Value = -3257.5314
$ws.Range("N137").Value = -13916.454
$ws.Range("H138").Value = 2311.2778
$ws.Range("I138").Value = 1890.5294
$ws.Range("J138").Value = 2687.7368
$ws.Range("K138").Value = 5671.5882
$ws.Range("L138").Value = 8063.2104
$ws.Range("M138").Value = -531.5882000000001
$ws.Range("N138").Value = -18343.2104

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 125004820
$ws.Range("I61").Value = 333341340
$ws.Range("J61").Value = 2898.2
$ws.Range("K61").Value = 333341340
$ws.Range("L61").Value = 2898.2
$ws.Range("M61").Value = -333341128
$ws.Range("N61").Value = -3322.2
$ws.Range("H97").Value = 581.8570999999999
$ws.Range("I97").Value = 581.8570999999999
$ws.Range("K97").Value = 581.8570999999999
$ws.Range("M97").Value = -85.85709999999995
$ws.Range("H136").Value = 125004820
$ws.Range("I136").Value = 333341340
$ws.Range("J136").Value = 2898.2
$ws.Range("K136").Value = 1000024020
$ws.Range("L136").Value = 8694.599999999999
$ws.Range("M136").Value = -1000021470
$ws.Range("N136").Value = -13794.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 46593.035
$ws.Range("I94").Value = 9846.348
$ws.Range("J94").Value = 257886.5
$ws.Range("K94").Value = 9846.348
$ws.Range("L94").Value = 257886.5
$ws.Range("M94").Value = -9395.348
$ws.Range("N94").Value = -258788.5
$ws.Range("H99").Value = 1960.4
$ws.Range("I99").Value = 1656.5714
$ws.Range("K99").Value = 1656.5714
$ws.Range("M99").Value = -158.5714
$ws.Range("H103").Value = 19400
$ws.Range("J103").Value = 19400
$ws.Range("L103").Value = 19400
$ws.Range("N103").Value = -21744
$ws.Range("H105").Value = 2513.1538
$ws.Range("I105").Value = 2513.1538
$ws.Range("K105").Value = 2513.1538
$ws.Range("M105").Value = -766.1538
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 256.1
$ws.Range("J7").Value = 580.25
$ws.Range("L7").Value = 580.25
$ws.Range("N7").Value = -806.25
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H99").Value = 4099.4
$ws.Range("I99").Value = 3999.25
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 3999.25
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -2501.25
$ws.Range("N99").Value = -7496
$ws.Range("H105").Value = 1819811.1
$ws.Range("I105").Value = 2501157
$ws.Range("K105").Value = 2501157
$ws.Range("M105").Value = -2499410
$ws.Range("H126").Value = 4099.4
$ws.Range("I126").Value = 3999.25
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 11997.75
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -9527.75
$ws.Range("N126").Value = -18440
$ws.Range("H134").Value = 7354765
$ws.Range("I134").Value = 9617277
$ws.Range("K134").Value = 28851831
$ws.Range("M134").Value = -28849296

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 92252.45
$ws.Range("I5").Value = 250578
$ws.Range("K5").Value = 751734
$ws.Range("M5").Value = -751622
$ws.Range("H33").Value = 2298.8
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H45").Value = 3683
$ws.Range("I45").Value = 1400
$ws.Range("J45").Value = 4009.1428
$ws.Range("K45").Value = 4200
$ws.Range("L45").Value = 12027.4284
$ws.Range("M45").Value = -3668
$ws.Range("N45").Value = -13091.4284
$ws.Range("H63").Value = 10000000
$ws.Range("I63").Value = 10000000
$ws.Range("K63").Value = 30000000
$ws.Range("M63").Value = -29999251
$ws.Range("H66").Value = 10000000
$ws.Range("I66").Value = 10000000
$ws.Range("K66").Value = 90000000
$ws.Range("M66").Value = -89996256
$ws.Range("H68").Value = 3599.75
$ws.Range("J68").Value = 2799.6667
$ws.Range("L68").Value = 8399.000100000001
$ws.Range("N68").Value = -10021.0001
$ws.Range("H71").Value = 3599.75
$ws.Range("J71").Value = 2799.6667
$ws.Range("L71").Value = 25197.0003
$ws.Range("N71").Value = -33309.0003
$ws.Range("H92").Value = 696.3333
$ws.Range("I92").Value = 399.5
$ws.Range("J92").Value = 1290
$ws.Range("K92").Value = 1198.5
$ws.Range("L92").Value = 3870
$ws.Range("M92").Value = 49.5
$ws.Range("N92").Value = -6366
$ws.Range("H97").Value = 1123.2858
$ws.Range("I97").Value = 662.3333
$ws.Range("K97").Value = 1986.9999
$ws.Range("M97").Value = -1490.9999
$ws.Range("H135").Value = 92252.45
$ws.Range("I135").Value = 250578
$ws.Range("K135").Value = 2255202
$ws.Range("M135").Value = -2252667

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2130.7334
$ws.Range("I107").Value = 2265.7273
$ws.Range("K107").Value = 2265.7273
$ws.Range("M107").Value = -345.7273
$ws.Range("H124").Value = 54999
$ws.Range("J124").Value = 54999
$ws.Range("L124").Value = 54999
$ws.Range("N124").Value = -64819
$ws.Range("H126").Value = 4653.5
$ws.Range("I126").Value = 4831.385
$ws.Range("K126").Value = 14494.155
$ws.Range("M126").Value = -12024.155

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H22").Value = 3634.25
$ws.Range("I22").Value = 3582
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 3582
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -3287
$ws.Range("N22").Value = -4590
$ws.Range("H27").Value = 3634.25
$ws.Range("I27").Value = 3582
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 3582
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = -3475
$ws.Range("N27").Value = -4214
$ws.Range("H46").Value = 1027.6428
$ws.Range("I46").Value = 1147.3334
$ws.Range("K46").Value = 1147.3334
$ws.Range("M46").Value = -959.3334
$ws.Range("H55").Value = 434.85715
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H68").Value = 2892252.5
$ws.Range("I68").Value = 3575481.8
$ws.Range("J68").Value = 500949.5
$ws.Range("K68").Value = 3575481.8
$ws.Range("L68").Value = 500949.5
$ws.Range("M68").Value = -3574732.8
$ws.Range("N68").Value = -502447.5
$ws.Range("H71").Value = 2892252.5
$ws.Range("I71").Value = 3575481.8
$ws.Range("J71").Value = 500949.5
$ws.Range("K71").Value = 17877409
$ws.Range("L71").Value = 2504747.5
$ws.Range("M71").Value = -17873665
$ws.Range("N71").Value = -2512235.5
$ws.Range("H93").Value = 2264.7778
$ws.Range("I93").Value = 2299.125
$ws.Range("J93").Value = 1990
$ws.Range("K93").Value = 2299.125
$ws.Range("L93").Value = 1990
$ws.Range("M93").Value = -1051.125
$ws.Range("N93").Value = -4486
$ws.Range("H100").Value = 12477935
$ws.Range("I100").Value = 16635922
$ws.Range("J100").Value = 3975
$ws.Range("K100").Value = 16635922
$ws.Range("L100").Value = 3975
$ws.Range("M100").Value = -16635381
$ws.Range("N100").Value = -5057
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 24689
$ws.Range("I26").Value = 24689
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 24689
$ws.Range("N26").ClearContents()
$ws.Range("M26").Value = -24396
$ws.Range("H81").Value = 145471.86
$ws.Range("I81").Value = 169233.5
$ws.Range("K81").Value = 338467
$ws.Range("M81").Value = -337406
$ws.Range("H84").Value = 145471.86
$ws.Range("I84").Value = 169233.5
$ws.Range("K84").Value = 1692335
$ws.Range("M84").Value = -1687031
$ws.Range("H125").Value = 201905
$ws.Range("J125").Value = 201905
$ws.Range("L125").Value = 201905
$ws.Range("N125").Value = -211745
$ws.Range("H132").Value = 17863760
$ws.Range("I132").Value = 26318986
$ws.Range("K132").Value = 78956958
$ws.Range("M132").Value = -78954428
$ws.Range("H136").Value = 15627867
$ws.Range("I136").Value = 16131911
$ws.Range("K136").Value = 48395733
$ws.Range("M136").Value = -48393183
